$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values so the table can hold 2D data
$ws.Range("B4").Value = 2.65
$ws.Range("B5").Value = 2.65
$ws.Range("C5").Value = 2.65

# Update selected cell from D15 to D5
$ws.Range("D5").Select()
